$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the sample-rate column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary rows 14-17: labels in column A, aggregate formulas in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format the summary values: bold, 12pt, vertically centered.
# Build the formatting once on B14, then fan it out with a format-only
# paste so every cell lands on the same shared style record.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the summary block
$ws.Range("A14:B17").RowHeight = 15.6

# Leave the summary block selected, matching the saved selection state
[void]$ws.Range("A14:B17").Select()

# Printable page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
